$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Library_Formula")

# Add new row 41: BR indicator override row (INDICATOR_12_OVER)
$a = $ws.Range("A41")
$a.Value = "CREATE/MODIFY"
$a.Font.Name = "Trebuchet MS"
$a.Font.Size = 10

$b = $ws.Range("B41")
$b.Value = "LIB_EWS_IT"
$b.Font.Name = "Trebuchet MS"
$b.Font.Size = 10

$c = $ws.Range("C41")
$c.Value = "INDICATOR_12_OVER"
$c.Font.Name = "Trebuchet MS"
$c.Font.Size = 10

$e = $ws.Range("E41")
$e.Value = "String"
$e.Font.Name = "Trebuchet MS"
$e.Font.Size = 10

$f = $ws.Range("F41")
$f.Value = "String"
$f.Font.Name = "Trebuchet MS"
$f.Font.Size = 10

# Move selection to the new row, matching the saved view state
$ws.Activate()
$ws.Range("D41").Select()
